# working on main div for content
#
# Reconstructs a small kanban-style task board: fills in the task names
# (col A), marks which status column (To Do / Doing / Done) each task is
# currently in with a checkmark, restyles those marker cells (centered
# alignment) and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: task names (rows 2-9) -----------------------------------
# Inserted in row order so the new shared-string table fills up in the
# same sequence as the source edit.
$ws.Range("A2").Value = "1. Create Header component"
$ws.Range("A3").Value = "2. Create Main component"
$ws.Range("A4").Value = "3. Create Footer component"
$ws.Range("A5").Value = "4. Create navbar component"
$ws.Range("A6").Value = "5. Create Buttons component "
$ws.Range("A7").Value = "6. Create Heading component"
$ws.Range("A8").Value = "7. Create Icons component"
$ws.Range("A9").Value = "8. Create Images component"

# --- Phase 2: first status checkmark (introduces the shared string) ---
$ws.Range("D2").Value = "✅"

# --- Phase 3: final (still unstarted) task row -------------------------
# A leading apostrophe forces this to be stored as text instead of being
# auto-coerced to the number 9 (Excel would otherwise treat "9. " as a
# numeric entry and drop the trailing period/space). Re-pasting the
# formatting from the row above afterwards removes the "stored as text"
# quote-prefix marker so the cell's style matches its neighbours again.
$ws.Range("A10").Value = "'9. "
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# --- Phase 4: remaining status checkmarks (reuse the shared string) ---
$ws.Range("D3").Value = "✅"
$ws.Range("B4").Value = "✅"
$ws.Range("C5").Value = "✅"
$ws.Range("C6").Value = "✅"
$ws.Range("C7").Value = "✅"
$ws.Range("B8").Value = "✅"
$ws.Range("B9").Value = "✅"

# --- Phase 5: restyle the marker cells ---------------------------------
# Build each combined style once on a scratch cell (outside the used
# range) and paste just the formatting onto every target cell, so all
# cells that need the same final style land on exactly one shared xf
# record instead of each other creating their own.

# Style used by B2 and D2:D5 (fill/border copied from B2, centered
# horizontally and vertically)
$ws.Range("B2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D2:D5").PasteSpecial(-4122)

# Style used by B3:B10 (fill/border copied from B3, centered
# horizontally and vertically)
$ws.Range("B3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").Copy()
$ws.Range("B3:B10").PasteSpecial(-4122)

# Style used by C2:C9 (fill/border copied from C2, centered
# horizontally, top-aligned vertically)
$ws.Range("C2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4160
$ws.Range("Z1").Copy()
$ws.Range("C2:C9").PasteSpecial(-4122)

$ws.Range("Z1").Clear()

# --- Phase 6: move the active selection --------------------------------
$ws.Range("B7").Select()

Write-Host "applied task board formatting"
